$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -8.340999999999999
$ws.Range("B4").Value = 8.656000000000001
$ws.Range("D4").Value = -7.895999999999999
$ws.Range("C5").Value = -12.828
$ws.Range("D5").Value = -7.382
$ws.Range("A8").Value = -21.624
$ws.Range("C8").Value = -12.837
$ws.Range("A10").Value = -21.705
$ws.Range("B11").Value = 7.114
$ws.Range("A12").Value = -21.382
$ws.Range("B12").Value = 6.334000000000001
$ws.Range("C12").Value = -11.495
$ws.Range("C13").Value = -12.518
$ws.Range("B15").Value = 5.274999999999999
$ws.Range("C15").Value = -12.966
$ws.Range("B17").Value = 4.809
$ws.Range("A18").Value = -21.744
$ws.Range("D20").Value = -7.901999999999999
$ws.Range("C21").Value = -12.973
$ws.Range("D23").Value = -8.265000000000001
$ws.Range("A25").Value = -21.779
$ws.Range("C25").Value = -12.181
$ws.Range("B26").Value = 6.025999999999999
$ws.Range("D26").Value = -7.797999999999999
$ws.Range("B27").Value = 5.194
$ws.Range("B28").Value = 4.789
$ws.Range("B32").Value = 6.586
$ws.Range("C32").Value = -12.207
$ws.Range("D34").Value = -7.853
$ws.Range("C36").Value = -13.173
$ws.Range("A37").Value = -21.492
$ws.Range("B37").Value = 5.474
$ws.Range("C38").Value = -11.859
$ws.Range("D39").Value = -7.476000000000001
$ws.Range("D40").Value = -8.01
$ws.Range("B41").Value = 9.232999999999999
$ws.Range("C41").Value = -12.982
$ws.Range("D41").Value = -8.047999999999998
$ws.Range("D42").Value = -8.056000000000001
$ws.Range("B47").Value = 5.757000000000001
$ws.Range("D47").Value = -8.071999999999999
$ws.Range("C50").Value = -13.026
$ws.Range("B51").Value = 5.913
$ws.Range("C52").Value = -11.553
$ws.Range("D52").Value = -7.197000000000001
$ws.Range("A55").Value = -22.008
$ws.Range("C59").Value = -12.571
$ws.Range("D60").Value = -8.125
$ws.Range("D62").Value = -8.051
$ws.Range("B65").Value = 5.922000000000001
$ws.Range("C67").Value = -10.944
$ws.Range("A68").Value = -21.576
$ws.Range("D70").Value = -7.358999999999999
$ws.Range("D72").Value = -7.607000000000001
$ws.Range("B73").Value = 7.016
$ws.Range("A77").Value = -20.915
$ws.Range("A78").Value = -20.644
$ws.Range("A79").Value = -20.994
$ws.Range("A80").Value = -20.938
$ws.Range("A81").Value = -21.811
$ws.Range("A82").Value = -21.705
$ws.Range("D83").Value = -8.131
$ws.Range("A84").Value = -21.703
$ws.Range("B84").Value = 5.994
$ws.Range("C84").Value = -12.371
$ws.Range("B85").Value = 4.849
$ws.Range("C86").Value = -13.31
$ws.Range("C88").Value = -13.371
$ws.Range("B89").Value = 4.55
$ws.Range("C89").Value = -13.317
$ws.Range("B93").Value = 5.728
$ws.Range("B95").Value = 6.571000000000001
$ws.Range("C95").Value = -12.407
$ws.Range("B98").Value = 7.261
$ws.Range("B99").Value = 5.709999999999999
$ws.Range("A101").Value = -21.897
$ws.Range("B101").Value = 5.046
$ws.Range("A102").Value = -20.661
$ws.Range("B102").Value = 7.733000000000001
$ws.Range("C105").Value = -12.653
